$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (prices) stay as text, matching the
# source data which stores every Price/Volume cell as a string.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '54.290.77'
$ws.Range("E2").Value = '  -0.03%  '
$ws.Range("D3").Value = '2.280.29'
$ws.Range("E3").Value = '  -0.20%  '
$ws.Range("E4").Value = '  +0.49%  '
$ws.Range("D5").Value = '502.53'
$ws.Range("E5").Value = '  +1.43%  '
$ws.Range("D6").Value = '129.28'
$ws.Range("E6").Value = '  +1.44%  '
$ws.Range("E7").Value = '  +0.29%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").Value = '0.0955'
$ws.Range("E9").Value = '  +1.53%  '
$ws.Range("D10").Value = '0.151'
$ws.Range("E10").Value = '  +1.03%  '
$ws.Range("E11").Value = '  +3.51%  '
$ws.Range("D12").Value = '4.74'
$ws.Range("E12").Value = '  +1.83%  '
$ws.Range("D13").Value = '2.690.21'
$ws.Range("E13").Value = '  +0.57%  '
$ws.Range("D14").Value = '22.95'
$ws.Range("E14").Value = '  +5.91%  '
$ws.Range("D15").Value = '54.254.48'
$ws.Range("E15").Value = '  +0.19%  '
$ws.Range("E16").Value = '  +0.32%  '
$ws.Range("D17").Value = '2.283.12'
$ws.Range("E17").Value = '  -0.61%  '
$ws.Range("D18").Value = '10.26'
$ws.Range("E18").Value = '  +3.46%  '
$ws.Range("E19").Value = '  +2.11%  '
$ws.Range("D20").Value = '304.70'
$ws.Range("E20").Value = '  +0.55%  '
$ws.Range("D21").Value = '6.41'
$ws.Range("E21").Value = '  +1.65%  '
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  +0.23%  '
$ws.Range("D23").Value = '62.04'
$ws.Range("E23").Value = '  -2.92%  '
$ws.Range("E24").Value = '  -0.21%  '
$ws.Range("E25").Value = '  +2.13%  '
$ws.Range("D26").Value = '7.33'
$ws.Range("E26").Value = '  +2.84%  '
$ws.Range("D27").Value = '174.99'
$ws.Range("E27").Value = '  +6.55%  '
$ws.Range("E28").Value = '  +1.48%  '
$ws.Range("D29").Value = '6.00'
$ws.Range("E29").Value = '  +2.06%  '
$ws.Range("D30").Value = '0.0₃0691'
$ws.Range("E30").Value = '  +1.68%  '
$ws.Range("E31").Value = '  +0.97%  '
$ws.Range("E32").Value = '  +0.10%  '
$ws.Range("D33").Value = '17.79'
$ws.Range("E33").Value = '  +1.58%  '
$ws.Range("E34").Value = '  -0.19%  '
$ws.Range("D35").Value = '0.942'
$ws.Range("E35").Value = '  +8.05%  '
$ws.Range("E36").Value = '  +1.17%  '
$ws.Range("D37").Value = '3.76'
$ws.Range("E37").Value = '  +3.28%  '
$ws.Range("D38").Value = '0.374'
$ws.Range("E38").Value = '  -0.57%  '
$ws.Range("E39").Value = '  +0.75%  '
$ws.Range("D40").Value = '3.39'
$ws.Range("E40").Value = '  +1.69%  '
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").Value = '125.19'
$ws.Range("E41").Value = '  -0.84%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = '4.81'
$ws.Range("E42").Value = '  -2.19%  '
$ws.Range("E43").Value = '  +3.30%  '
$ws.Range("D44").Value = '0.0895'
$ws.Range("E44").Value = '  +0.57%  '
$ws.Range("D45").Value = '0.548'
$ws.Range("E45").Value = '  +0.33%  '
$ws.Range("D46").Value = '240.26'
$ws.Range("E46").Value = '  +0.98%  '
$ws.Range("D47").Value = '0.372'
$ws.Range("E47").Value = '  -0.71%  '
$ws.Range("D48").Value = '0.0206'
$ws.Range("E48").Value = '  +0.99%  '
$ws.Range("D49").Value = '10.78'
$ws.Range("E49").Value = '  +1.00%  '
$ws.Range("D50").Value = '16.39'
$ws.Range("E50").Value = '  +0.14%  '
$ws.Range("D51").Value = '4.65'
$ws.Range("E51").Value = '  +0.58%  '
